$wb = $excel.ActiveWorkbook

# --- Update the credibility_renovation_strategy sheet ---
# The "submitted / not yet assessed" status is retired: every country that had
# it now shows "submitted / assessed". Additionally HU, MT and SI move from
# "not submitted" to "submitted / assessed" (PL stays "not submitted").
$wsCred = $wb.Worksheets.Item("credibility_renovation_strategy")

$rowsToUpdate = @(2, 3, 8, 11, 13, 15, 16, 18, 19, 23, 24, 25, 26)
foreach ($r in $rowsToUpdate) {
    $wsCred.Range("B$r").Value = "submitted / assessed"
}

# --- Reflect the final active sheet/selection state ---
$wsFossil = $wb.Worksheets.Item("limit_fossil_fuels_buildings")
$wsFossil.Range("G6").Select()

$wsCred.Activate()
$wsCred.Range("B31").Select()
